$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "19×73="; New = "49×14=" },
    @{ Old = "17×66="; New = "75×98=" },
    @{ Old = "83×67="; New = "23×25=" },
    @{ Old = "27×98="; New = "65×29=" },
    @{ Old = "41×27="; New = "65×62=" },
    @{ Old = "80×75="; New = "16×61=" },
    @{ Old = "31×20="; New = "75×93=" },
    @{ Old = "11×79="; New = "16×37=" },
    @{ Old = "55×28="; New = "91×89=" },
    @{ Old = "84×48="; New = "18×72=" },
    @{ Old = "45×68="; New = "29×58=" },
    @{ Old = "50×98="; New = "38×26=" },
    @{ Old = "63×92="; New = "69×64=" },
    @{ Old = "28×44="; New = "83×56=" },
    @{ Old = "64×40="; New = "93×56=" },
    @{ Old = "78×69="; New = "73×23=" },
    @{ Old = "62×51="; New = "39×93=" },
    @{ Old = "30×75="; New = "84×49=" },
    @{ Old = "19×45="; New = "76×45=" },
    @{ Old = "67×81="; New = "24×65=" },
    @{ Old = "25×41="; New = "99×21=" },
    @{ Old = "66×64="; New = "80×74=" },
    @{ Old = "74×50="; New = "51×53=" },
    @{ Old = "44×56="; New = "17×22=" },
    @{ Old = "34×20="; New = "94×78=" }
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.New, 2)
}
